{"js": "// Update each lattice-multiplication exercise cell in-place (same 5x3 table\n// shape/formatting, just new operands/digits per the commit's regenerated\n// exercise set). Cell text uses \"\\v\" (vertical-tab / Word line-break char)\n// to join the lines that are separated by <w:br/> inside the cell's run.\nconst cellsData = [\n  { row: 0, col: 0, lines: [\"23 x 31\", \"  3    1\", \"  ----\", \"2|    |\", \"3|    |\"] },\n  { row: 0, col: 1, lines: [\"60 x 80\", \"  8    0\", \"  ----\", \"6|    |\", \"0|    |\"] },\n  { row: 0, col: 2, lines: [\"40 x 74\", \"  7    4\", \"  ----\", \"4|    |\", \"0|    |\"] },\n  { row: 1, col: 0, lines: [\"44 x 76\", \"  7    6\", \"  ----\", \"4|    |\", \"4|    |\"] },\n  { row: 1, col: 1, lines: [\"74 x 64\", \"  6    4\", \"  ----\", \"7|    |\", \"4|    |\"] },\n  { row: 1, col: 2, lines: [\"45 x 95\", \"  9    5\", \"  ----\", \"4|    |\", \"5|    |\"] },\n  { row: 2, col: 0, lines: [\"44 x 19\", \"  1    9\", \"  ----\", \"4|    |\", \"4|    |\"] },\n  { row: 2, col: 1, lines: [\"54 x 94\", \"  9    4\", \"  ----\", \"5|    |\", \"4|    |\"] },\n  { row: 2, col: 2, lines: [\"73 x 63\", \"  6    3\", \"  ----\", \"7|    |\", \"3|    |\"] },\n  { row: 3, col: 0, lines: [\"20 x 40\", \"  4    0\", \"  ----\", \"2|    |\", \"0|    |\"] },\n  { row: 3, col: 1, lines: [\"88 x 83\", \"  8    3\", \"  ----\", \"8|    |\", \"8|    |\"] },\n  { row: 3, col: 2, lines: [\"53 x 97\", \"  9    7\", \"  ----\", \"5|    |\", \"3|    |\"] },\n  { row: 4, col: 0, lines: [\"20 x 76\", \"  7    6\", \"  ----\", \"2|    |\", \"0|    |\"] },\n  { row: 4, col: 1, lines: [\"31 x 14\", \"  1    4\", \"  ----\", \"3|    |\", \"1|    |\"] },\n  { row: 4, col: 2, lines: [\"40 x 88\", \"  8    8\", \"  ----\", \"4|    |\", \"0|    |\"] },\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.body.paragraphs.load(\"items\"); // no-op load just to be safe on some hosts\n\nfor (const { row, col, lines } of cellsData) {\n  const cell = table.getCell(row, col);\n  cell.body.paragraphs.load(\"items\");\n  await context.sync();\n\n  const para = cell.body.paragraphs.items[0];\n  const range = para.getRange();\n  range.insertText(lines.join(\"\\u000b\"), Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update each lattice-multiplication exercise cell in-place (same 5x3 table\n# shape/formatting, just new operands/digits per the commit's regenerated\n# exercise set). [char]11 is the vertical-tab / Word manual line-break\n# character used between the lines that <w:br/> separates inside a cell.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$br = [char]11\n\n$cellsData = @(\n    @{ Row = 1; Col = 1; Lines = @(\"23 x 31\", \"  3    1\", \"  ----\", \"2|    |\", \"3|    |\") },\n    @{ Row = 1; Col = 2; Lines = @(\"60 x 80\", \"  8    0\", \"  ----\", \"6|    |\", \"0|    |\") },\n    @{ Row = 1; Col = 3; Lines = @(\"40 x 74\", \"  7    4\", \"  ----\", \"4|    |\", \"0|    |\") },\n    @{ Row = 2; Col = 1; Lines = @(\"44 x 76\", \"  7    6\", \"  ----\", \"4|    |\", \"4|    |\") },\n    @{ Row = 2; Col = 2; Lines = @(\"74 x 64\", \"  6    4\", \"  ----\", \"7|    |\", \"4|    |\") },\n    @{ Row = 2; Col = 3; Lines = @(\"45 x 95\", \"  9    5\", \"  ----\", \"4|    |\", \"5|    |\") },\n    @{ Row = 3; Col = 1; Lines = @(\"44 x 19\", \"  1    9\", \"  ----\", \"4|    |\", \"4|    |\") },\n    @{ Row = 3; Col = 2; Lines = @(\"54 x 94\", \"  9    4\", \"  ----\", \"5|    |\", \"4|    |\") },\n    @{ Row = 3; Col = 3; Lines = @(\"73 x 63\", \"  6    3\", \"  ----\", \"7|    |\", \"3|    |\") },\n    @{ Row = 4; Col = 1; Lines = @(\"20 x 40\", \"  4    0\", \"  ----\", \"2|    |\", \"0|    |\") },\n    @{ Row = 4; Col = 2; Lines = @(\"88 x 83\", \"  8    3\", \"  ----\", \"8|    |\", \"8|    |\") },\n    @{ Row = 4; Col = 3; Lines = @(\"53 x 97\", \"  9    7\", \"  ----\", \"5|    |\", \"3|    |\") },\n    @{ Row = 5; Col = 1; Lines = @(\"20 x 76\", \"  7    6\", \"  ----\", \"2|    |\", \"0|    |\") },\n    @{ Row = 5; Col = 2; Lines = @(\"31 x 14\", \"  1    4\", \"  ----\", \"3|    |\", \"1|    |\") },\n    @{ Row = 5; Col = 3; Lines = @(\"40 x 88\", \"  8    8\", \"  ----\", \"4|    |\", \"0|    |\") }\n)\n\nforeach ($entry in $cellsData) {\n    $cell = $t.Cell($entry.Row, $entry.Col)\n    $cell.Range.Text = [string]::Join($br, $entry.Lines)\n}\n"}
